$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2..27, columns: F, H, I, K, L, M, N, O, P, Q, R, S
$data = @(
    @(202, 'bedrooms', 'target', 'j', 'stimuli/img_ybbmx.png', 55.24324324324324, 36.75675675675676, 46, 37, 3, 3, 3),
    @(203, 'bedrooms', 'target', 'j', 'stimuli/img_h0hbk.png', 86.80952380952381, 69.19047619047619, 78, 42, 9, 9, 9),
    @(204, 'bedrooms', 'target', 'j', 'stimuli/img_5p2ql.png', 89.19565217391305, 72.52173913043478, 80.8586956521739, 46, 10, 10, 10),
    @(205, 'bedrooms', 'target', 'j', 'stimuli/img_2pk6v.png', 85.08108108108108, 66.16216216216216, 75.62162162162161, 37, 9, 9, 9),
    @(206, 'bedrooms', 'target', 'j', 'stimuli/img_okvvw.png', 50.58333333333334, 32.11111111111111, 41.34722222222223, 36, 2, 2, 2),
    @(207, 'bedrooms', 'target', 'j', 'stimuli/img_x0u5z.png', 92, 78.16216216216216, 85.08108108108108, 37, 10, 10, 10),
    @(208, 'bedrooms', 'target', 'j', 'stimuli/img_rvssl.png', 74.25, 54.33333333333334, 64.29166666666667, 36, 6, 6, 6),
    @(209, 'bedrooms', 'target', 'j', 'stimuli/img_uxxo0.png', 71.74418604651163, 48.44186046511628, 60.09302325581395, 43, 5, 5, 5),
    @(210, 'bedrooms', 'target', 'j', 'stimuli/img_71mhq.png', 69.34210526315789, 47.02631578947368, 58.18421052631579, 38, 5, 5, 5),
    @(211, 'bedrooms', 'target', 'j', 'stimuli/img_qgbyn.png', 65.08108108108108, 40.10810810810811, 52.5945945945946, 37, 4, 4, 4),
    @(212, 'kitchens', 'distractor', 'f', 'stimuli/img_p3hpc.png', 72.83333333333333, 52.22222222222222, 62.52777777777777, 36, 6, 6, 6),
    @(213, 'bedrooms', 'target', 'j', 'stimuli/img_v8dra.png', 61.77272727272727, 38.79545454545455, 50.28409090909091, 44, 3, 3, 3),
    @(214, 'bedrooms', 'target', 'j', 'stimuli/img_2js6m.png', 40.02777777777778, 20.88888888888889, 30.45833333333334, 36, 2, 2, 2),
    @(215, 'bedrooms', 'target', 'j', 'stimuli/img_fqgem.png', 80.75, 61.475, 71.1125, 40, 8, 8, 8),
    @(216, 'bedrooms', 'target', 'j', 'stimuli/img_le8uf.png', 12.88888888888889, 9.222222222222221, 11.05555555555556, 36, 1, 1, 1),
    @(217, 'bedrooms', 'target', 'j', 'stimuli/img_wyctg.png', 33.44736842105263, 11.39473684210526, 22.42105263157895, 38, 1, 1, 1),
    @(218, 'bedrooms', 'target', 'j', 'stimuli/img_th7xh.png', 82.35897435897436, 65.53846153846153, 73.94871794871796, 39, 8, 8, 8),
    @(219, 'bedrooms', 'target', 'j', 'stimuli/img_a9acb.png', 77.11428571428571, 58.42857142857143, 67.77142857142857, 35, 7, 7, 7),
    @(220, 'kitchens', 'distractor', 'f', 'stimuli/img_di49w.png', 73.03125, 52.25, 62.640625, 32, 6, 6, 6),
    @(221, 'living_rooms', 'distractor', 'f', 'stimuli/img_x4bln.png', 76.34042553191489, 59.51063829787234, 67.92553191489361, 47, 7, 7, 7),
    @(222, 'living_rooms', 'distractor', 'f', 'stimuli/img_z4jxm.png', 88.30952380952381, 72.64285714285714, 80.47619047619048, 42, 10, 10, 10),
    @(223, 'bedrooms', 'target', 'j', 'stimuli/img_bj2gr.png', 65.25, 44.8, 55.025, 40, 4, 4, 4),
    @(224, 'kitchens', 'distractor', 'f', 'stimuli/img_2b8fp.png', 73.89189189189189, 51.45945945945946, 62.67567567567568, 37, 6, 6, 6),
    @(225, 'bedrooms', 'target', 'j', 'stimuli/img_oou46.png', 75.70270270270271, 54.86486486486486, 65.28378378378379, 37, 6, 6, 6),
    @(226, 'bedrooms', 'target', 'j', 'stimuli/img_5m6x4.png', 80.23076923076923, 58.41025641025641, 69.32051282051282, 39, 7, 7, 7),
    @(227, 'kitchens', 'distractor', 'f', 'stimuli/img_cnyac.png', 69.1470588235294, 47.8235294117647, 58.48529411764706, 34, 5, 5, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 6).Value = $row[0]    # F trial_total
    $ws.Cells.Item($r, 8).Value = $row[1]    # H category
    $ws.Cells.Item($r, 9).Value = $row[2]    # I cond_cat
    $ws.Cells.Item($r, 11).Value = $row[3]   # K correct_answer
    $ws.Cells.Item($r, 12).Value = $row[4]   # L stimulus
    $ws.Cells.Item($r, 13).Value = $row[5]   # M conceptual
    $ws.Cells.Item($r, 14).Value = $row[6]   # N perceptual
    $ws.Cells.Item($r, 15).Value = $row[7]   # O typicality
    $ws.Cells.Item($r, 16).Value = $row[8]   # P n
    $ws.Cells.Item($r, 17).Value = $row[9]   # Q p_typicality
    $ws.Cells.Item($r, 18).Value = $row[10]  # R p_conceptual
    $ws.Cells.Item($r, 19).Value = $row[11]  # S p_perceptual
}
